$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.154.31"
$ws.Range("E2").Value = "  -0.99%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.861.04"
$ws.Range("E3").Value = "  -0.71%  "

# Row 4
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7058"
$ws.Range("E5").Value = "  -0.98%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "241.15"
$ws.Range("E6").Value = "  -0.57%  "

# Row 7
$ws.Range("E7").Value = "  -0.05%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3079"
$ws.Range("E8").Value = "  -1.15%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07626"
$ws.Range("E9").Value = "  -3.62%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.62"
$ws.Range("E10").Value = "  -1.63%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08340"
$ws.Range("E11").Value = "  +0.91%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.877.85"
$ws.Range("E12").Value = "  +0.15%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.177"
$ws.Range("E13").Value = "  -2.14%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7075"
$ws.Range("E14").Value = "  -2.67%  "

# Row 15
$ws.Range("E15").Value = "  +0.11%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.178.13"
$ws.Range("E16").Value = "  -0.97%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.903"
$ws.Range("E17").Value = "  -0.49%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "242.54"
$ws.Range("E18").Value = "  -2.19%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007798"
$ws.Range("E19").Value = "  -1.00%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.115.00"
$ws.Range("E20").Value = "  -0.83%  "

# Row 21
$ws.Range("E21").Value = "  -1.77%  "

# Row 22
$ws.Range("E22").Value = "  +0.00%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.845"
$ws.Range("E23").Value = "  -1.30%  "

# Row 24
$ws.Range("E24").Value = "  +0.01%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1587"
$ws.Range("E25").Value = "  +0.05%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.99"
$ws.Range("E26").Value = "  -0.64%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.914"
$ws.Range("E27").Value = "  -1.13%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.44"
$ws.Range("E28").Value = "  +0.55%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.337"
$ws.Range("E29").Value = "  -1.93%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.499"
$ws.Range("E30").Value = "  +0.20%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.393"
$ws.Range("E31").Value = "  +0.22%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.203"
$ws.Range("E32").Value = "  +1.99%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05130"
$ws.Range("E33").Value = "  -3.43%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7937"
$ws.Range("E34").Value = "  +9.79%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.907"
$ws.Range("E35").Value = "  -1.33%  "

# Row 36
$ws.Range("E36").Value = "  -3.25%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.681"
$ws.Range("E37").Value = "  -0.09%  "

# Row 38
$ws.Range("E38").Value = "  -1.77%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.695"
$ws.Range("E39").Value = "  -1.03%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.165.40"
$ws.Range("E40").Value = "  -6.05%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.181"
$ws.Range("E41").Value = "  +0.08%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8897"
$ws.Range("E42").Value = "  -2.31%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "72.71"
$ws.Range("E43").Value = "  -1.75%  "

# Row 44
$ws.Range("E44").Value = "  -0.12%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.82"
$ws.Range("E45").Value = "  -1.26%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.010.97"
$ws.Range("E46").Value = "  -1.00%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5192"
$ws.Range("E47").Value = "  -2.59%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.763"
$ws.Range("E48").Value = "  -0.08%  "

# Row 49
$ws.Range("E49").Value = "  +0.27%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.326"
$ws.Range("E50").Value = "  +0.28%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.000"
$ws.Range("E51").Value = "  -0.34%  "
